$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking values stored as TEXT in the source
# data (e.g. "98.506.05", "1.00", "0.0000246"). Force each touched Price cell to
# Text format before assigning so Excel does not reinterpret it as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "98.506.05"
$ws.Range("E2").Value = "  +5.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.356.85"
$ws.Range("E3").Value = "  +10.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.49"
$ws.Range("E5").Value = "  +9.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.12"
$ws.Range("E6").Value = "  +3.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.19"
$ws.Range("E7").Value = "  +10.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.386"
$ws.Range("E8").Value = "  +4.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.353.33"
$ws.Range("E10").Value = "  +10.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.807"
$ws.Range("E11").Value = "  +1.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.200"
$ws.Range("E12").Value = "  +2.96%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "98.211.48"
$ws.Range("E13").Value = "  +5.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.85"
$ws.Range("E14").Value = "  +8.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000246"
$ws.Range("E15").Value = "  +4.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.975.19"
$ws.Range("E16").Value = "  +10.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.50"
$ws.Range("E17").Value = "  +5.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.358.49"
$ws.Range("E18").Value = "  +10.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.58"
$ws.Range("E19").Value = "  +3.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.00"
$ws.Range("E20").Value = "  +5.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "485.46"
$ws.Range("E21").Value = "  +11.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.85"
$ws.Range("E22").Value = "  +3.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000207"
$ws.Range("E23").Value = "  +11.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("E24").Value = "  +5.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.67"
$ws.Range("E25").Value = "  +4.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.22"
$ws.Range("E26").Value = "  +5.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.00"
$ws.Range("E27").Value = "  +3.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.530.11"
$ws.Range("E28").Value = "  +10.44%  "

$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.252"
$ws.Range("E30").Value = "  +3.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("B31").Value = "Cronos"
$ws.Range("C31").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D31").Value = "0.187"
$ws.Range("E31").Value = "  +7.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.126"
$ws.Range("E32").Value = "  +4.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").Value = "  -11.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.28"
$ws.Range("E34").Value = "  +4.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.30"
$ws.Range("E35").Value = "  +9.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "518.78"
$ws.Range("E36").Value = "  +13.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "7.39"
$ws.Range("E37").Value = "  -1.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.152"
$ws.Range("E38").Value = "  -0.41%  "

$ws.Range("E39").Value = "  +3.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.88"
$ws.Range("E40").Value = "  +3.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.447"
$ws.Range("E41").Value = "  +4.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.26"
$ws.Range("E42").Value = "  +2.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.64"
$ws.Range("E43").Value = "  -3.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.24"
$ws.Range("E44").Value = "  +5.38%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.775"
$ws.Range("E46").Value = "  +18.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.29"
$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  +7.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "1.37"
$ws.Range("E49").Value = "  +9.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "45.49"
$ws.Range("E50").Value = "  +4.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.52"
$ws.Range("E51").Value = "  +8.29%  "
